$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename column header value B3 from "patient_id" to "participant_id"
$ws.Range("B3").Value = "participant_id"

# Update the active selection/view: select cell M4 (clears previous A8:XFD8
# selection and the scrolled topLeftCell position)
$ws.Range("M4").Select()
